$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of column A and column B for every data row (1-66).
# The correspondence table was flipped (orthography regex conversion direction reversed).
for ($r = 1; $r -le 66; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value()
    $valB = $cellB.Value()

    $cellA.Value = $valB
    $cellB.Value = $valA
}

# Row 62 carried a special font color (dark gray, RGB 333333) on column A;
# after the swap that formatting follows the value into column B.
$ws.Cells.Item(62, 1).Style = "Normal"
$ws.Cells.Item(62, 2).Font.Color = 3355443

# Row 63 carried a special font color (near-black, RGB 222222) on column B;
# after the swap that formatting follows the value into column A.
$ws.Cells.Item(63, 2).Style = "Normal"
$ws.Cells.Item(63, 1).Font.Color = 2236962

# The selected range moves from column C down to columns A:B.
$ws.Range("A1:B1048576").Select()

Write-Output "done"
